$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2..439) from 45190 to 45192
$ws.Range("C2:C439").Value = 45192

# Add the new row 440 with data
$ws.Cells.Item(440, 1).Value = "A 45218-2023"
$ws.Cells.Item(440, 2).Value = 45191
$ws.Cells.Item(440, 3).Value = 45192
$ws.Cells.Item(440, 4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item(440, 5).Value = "KALIX"
$ws.Cells.Item(440, 7).Value = 1.5
$ws.Cells.Item(440, 8).Value = 0
$ws.Cells.Item(440, 9).Value = 0
$ws.Cells.Item(440, 10).Value = 0
$ws.Cells.Item(440, 11).Value = 0
$ws.Cells.Item(440, 12).Value = 0
$ws.Cells.Item(440, 13).Value = 0
$ws.Cells.Item(440, 14).Value = 0
$ws.Cells.Item(440, 15).Value = 0
$ws.Cells.Item(440, 16).Value = 0
$ws.Cells.Item(440, 17).Value = 0

# Match date formatting (style index 1, numFmt yyyy-mm-dd) used by other rows in columns B and C
$ws.Cells.Item(440, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(440, 3).NumberFormat = "YYYY-MM-DD"

# Match the wrap-text style used for column R on data rows
$ws.Cells.Item(440, 18).WrapText = $true

# Normalize row height for row 439 (gains an explicit custom row height in the target file)
$ws.Rows.Item(439).RowHeight = 15
